$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (100kΩ Resistor / C149504): quantity changes from 1 to 2
$ws.Range("C7").Value = 2

# Row 9 previously held "Switch" / "C319028"; replace with the new part
$ws.Range("A9").Value = "39kΩ Resistor"
$ws.Range("B9").Value = "C25826"
$ws.Range("D9").Value = 0.0018

# Reflect the last-used selection
$ws.Range("D9").Select()
